# Fixes for VRelay and Upgrade Verification SCript
#
# This script updates the timestamp cells (columns B and D) across the
# Katalon "VT-Data-Prod" bootstrap sheets, simulating the latest run of the
# verification script which stamped new "last executed" dates into the
# workbook. Cells are written in the same chronological order the
# automation tool produced them so that the resulting shared-string table
# is built up the same way.

$wb = $excel.ActiveWorkbook

$shDualCF       = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
$shNoCF         = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
$shSingleCF     = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
$shCreditDualCF = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
$shCreditNoCF   = $wb.Worksheets.Item("VT-SaleCredit-NoCF-Generic")
$shCreditSingle = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
$shAuthCapCred  = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
$shAuthCapVoid  = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
$shManualAuth   = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")

# --- Pass 1 (Wed Aug 20 23:19 - 23:31 IST 2025): VT-AuthCapCredit-Generic and
# VT-AuthCapVoid-Generic "DateDemo" (column D) timestamps ---
$shAuthCapCred.Range("D2").Value = "Wed Aug 20 23:19:44 IST 2025"
$shAuthCapCred.Range("D3").Value = "Wed Aug 20 23:20:58 IST 2025"
$shAuthCapCred.Range("D4").Value = "Wed Aug 20 23:22:01 IST 2025"
$shAuthCapCred.Range("D5").Value = "Wed Aug 20 23:23:07 IST 2025"
$shAuthCapCred.Range("D6").Value = "Wed Aug 20 23:24:16 IST 2025"
$shAuthCapCred.Range("D7").Value = "Wed Aug 20 23:25:29 IST 2025"

$shAuthCapVoid.Range("D2").Value = "Wed Aug 20 23:26:32 IST 2025"
$shAuthCapVoid.Range("D3").Value = "Wed Aug 20 23:27:33 IST 2025"
$shAuthCapVoid.Range("D4").Value = "Wed Aug 20 23:28:31 IST 2025"
$shAuthCapVoid.Range("D5").Value = "Wed Aug 20 23:29:31 IST 2025"
$shAuthCapVoid.Range("D6").Value = "Wed Aug 20 23:30:36 IST 2025"
$shAuthCapVoid.Range("D7").Value = "Wed Aug 20 23:31:37 IST 2025"

# --- Pass 2 (Wed Aug 20 23:58 - Thu Aug 21 00:20 IST 2025): VT-SaleCredit-NoCF-Generic
# "DateProd" (column B) and VT-SaleVoid-NoCF-Generic "DateDemo" (column D) ---
$shCreditNoCF.Range("B2").Value = "Thu Aug 21 00:06:36 IST 2025"
$shCreditNoCF.Range("B3").Value = "Thu Aug 21 00:07:15 IST 2025"
$shCreditNoCF.Range("B4").Value = "Thu Aug 21 00:07:54 IST 2025"
$shCreditNoCF.Range("B5").Value = "Thu Aug 21 00:08:38 IST 2025"

$shNoCF.Range("D2").Value = "Thu Aug 21 00:14:57 IST 2025"
$shNoCF.Range("D3").Value = "Thu Aug 21 00:15:38 IST 2025"
$shNoCF.Range("D4").Value = "Thu Aug 21 00:16:17 IST 2025"
$shNoCF.Range("D5").Value = "Thu Aug 21 00:17:01 IST 2025"

# --- Pass 3 (Fri Aug 22 23:07 - 23:30 IST 2025): VT-AuthCapCredit-Generic and
# VT-AuthCapVoid-Generic "DateProd" (column B) timestamps ---
$shAuthCapCred.Range("B2").Value = "Fri Aug 22 23:15:32 IST 2025"
$shAuthCapCred.Range("B3").Value = "Fri Aug 22 23:16:59 IST 2025"
$shAuthCapCred.Range("B4").Value = "Fri Aug 22 23:18:40 IST 2025"
$shAuthCapCred.Range("B5").Value = "Fri Aug 22 23:19:56 IST 2025"
$shAuthCapCred.Range("B6").Value = "Fri Aug 22 23:21:23 IST 2025"
$shAuthCapCred.Range("B7").Value = "Fri Aug 22 23:22:47 IST 2025"

$shAuthCapVoid.Range("B2").Value = "Fri Aug 22 23:24:11 IST 2025"
$shAuthCapVoid.Range("B3").Value = "Fri Aug 22 23:25:31 IST 2025"
$shAuthCapVoid.Range("B4").Value = "Fri Aug 22 23:26:52 IST 2025"
$shAuthCapVoid.Range("B5").Value = "Fri Aug 22 23:28:14 IST 2025"
$shAuthCapVoid.Range("B6").Value = "Fri Aug 22 23:29:28 IST 2025"
$shAuthCapVoid.Range("B7").Value = "Fri Aug 22 23:30:59 IST 2025"

# --- Pass 4 (Sat Aug 23 00:04 - 00:32 IST 2025): remaining "DateProd"
# (column B) timestamps across the rest of the sheets ---
$shManualAuth.Range("B2").Value = "Sat Aug 23 00:04:18 IST 2025"
$shManualAuth.Range("B3").Value = "Sat Aug 23 00:05:28 IST 2025"
$shManualAuth.Range("B4").Value = "Sat Aug 23 00:06:25 IST 2025"
$shManualAuth.Range("B5").Value = "Sat Aug 23 00:07:21 IST 2025"
$shManualAuth.Range("B6").Value = "Sat Aug 23 00:08:30 IST 2025"
$shManualAuth.Range("B7").Value = "Sat Aug 23 00:09:22 IST 2025"

$shCreditDualCF.Range("B2").Value = "Sat Aug 23 00:10:20 IST 2025"
$shCreditDualCF.Range("B3").Value = "Sat Aug 23 00:11:25 IST 2025"
$shCreditDualCF.Range("B4").Value = "Sat Aug 23 00:12:16 IST 2025"
$shCreditDualCF.Range("B5").Value = "Sat Aug 23 00:13:09 IST 2025"

$shCreditSingle.Range("B2").Value = "Sat Aug 23 00:18:20 IST 2025"
$shCreditSingle.Range("B3").Value = "Sat Aug 23 00:19:24 IST 2025"
$shCreditSingle.Range("B4").Value = "Sat Aug 23 00:20:22 IST 2025"
$shCreditSingle.Range("B5").Value = "Sat Aug 23 00:21:18 IST 2025"

$shDualCF.Range("B2").Value = "Sat Aug 23 00:22:08 IST 2025"
$shDualCF.Range("B3").Value = "Sat Aug 23 00:23:17 IST 2025"
$shDualCF.Range("B4").Value = "Sat Aug 23 00:24:18 IST 2025"
$shDualCF.Range("B5").Value = "Sat Aug 23 00:25:23 IST 2025"

$shNoCF.Range("B2").Value = "Sat Aug 23 00:26:19 IST 2025"
$shNoCF.Range("B3").Value = "Sat Aug 23 00:27:05 IST 2025"
$shNoCF.Range("B4").Value = "Sat Aug 23 00:28:00 IST 2025"
$shNoCF.Range("B5").Value = "Sat Aug 23 00:28:46 IST 2025"

$shSingleCF.Range("B2").Value = "Sat Aug 23 00:29:30 IST 2025"
$shSingleCF.Range("B3").Value = "Sat Aug 23 00:30:21 IST 2025"
$shSingleCF.Range("B4").Value = "Sat Aug 23 00:31:25 IST 2025"
$shSingleCF.Range("B5").Value = "Sat Aug 23 00:32:14 IST 2025"
